$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row stays the same text (Ano/Brasil/Nordeste/Sergipe) in A1:D1 -
# re-assert values so the sheet is rewritten the way the source file was.
$ws.Range("A1").Value = "Ano"
$ws.Range("B1").Value = "Brasil"
$ws.Range("C1").Value = "Nordeste"
$ws.Range("D1").Value = "Sergipe"

# Updated figures for 2019 (row 9) and 2020 (row 10)
$ws.Range("B9").Value = 0.74652644359447518
$ws.Range("C9").Value = -1.4189382427712727
$ws.Range("D9").Value = -0.12419148437128325

$ws.Range("B10").Value = -4.3024599759648456
$ws.Range("C10").Value = -5.4472613016193216
$ws.Range("D10").Value = -10.013197113568694

# The wide B2:D12 selection is dropped in favour of a plain A1 selection
$ws.Range("A1").Select()

# Columns get explicit, fitted widths (A narrow "Ano" column, B-D sized to
# their formatted figures)
$ws.Columns("A").ColumnWidth = 4.166666666666667
$ws.Columns("B").ColumnWidth = 5.0
$ws.Columns("C").ColumnWidth = 8.5
$ws.Columns("D").ColumnWidth = 6.833333333333334
